# Add the new "2021年" row (row 12) to the sheet, matching the
# formatting of the preceding row (row 11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 11 down into the new row 12 so the new row's
# label cell (A12) picks up the same style (centered, bold, bordered) as
# the other year labels in column A.
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 49
$ws.Range("D12").Value = 1975
$ws.Range("G12").Value = 1926
